# Weekly driver report update for 2025-04-21
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# Bad Drivers section
$ws.Range("C5").Value = 2509
$ws.Range("D5").Value = 98.5
$ws.Range("C6").Value = 2514

# Good Drivers section (Total Samples column B)
$ws.Range("B16").Value = 56069
$ws.Range("B21").Value = 276086
$ws.Range("B22").Value = 625298
$ws.Range("B27").Value = 331283
$ws.Range("B29").Value = 453652
$ws.Range("B37").Value = 96091
$ws.Range("B40").Value = 99549
$ws.Range("B43").Value = 175767
$ws.Range("B44").Value = 240182
$ws.Range("B52").Value = 684728
$ws.Range("B54").Value = 210188
$ws.Range("B57").Value = 308481
$ws.Range("B64").Value = 443223
$ws.Range("B66").Value = 109665
$ws.Range("B68").Value = 62515
